# Auto-generated edit script applying the Brynhildr_Profits.xlsx diff
# Updates numeric leve-profit calculation columns (H-N) across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets: revised average/median market prices ripple
# into NQ/HQ price and profit columns; a few rows gain or lose a trailing cell
# where a profit value became newly computable (or no longer applicable).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 155.4
$ws.Range("H10").Value = 15000
$ws.Range("I10").Value = 15000
$ws.Range("K10").Value = 15000
$ws.Range("M10").Value = -14707
$ws.Range("H12").Value = 5011.75
$ws.Range("I12").Value = 182.33333
$ws.Range("K12").Value = 182.33333
$ws.Range("M12").Value = -12.33332999999999
$ws.Range("H19").Value = 1103.3334
$ws.Range("I19").Value = 605.7143
$ws.Range("K19").Value = 605.7143
$ws.Range("M19").Value = -430.7143
$ws.Range("H53").Value = 218.63637
$ws.Range("I53").Value = 199.375
$ws.Range("K53").Value = 199.375
$ws.Range("M53").Value = 437.625
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H70").Value = 3200
$ws.Range("I70").Value = 2050.3333
$ws.Range("K70").Value = 6150.999899999999
$ws.Range("M70").Value = -5880.999899999999
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H73").Value = 3200
$ws.Range("I73").Value = 2050.3333
$ws.Range("K73").Value = 6150.999899999999
$ws.Range("M73").Value = -5214.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 149112.39
$ws.Range("I32").Value = 158827.58
$ws.Range("K32").Value = 158827.58
$ws.Range("M32").Value = -158540.58
$ws.Range("H61").Value = 7442.524
$ws.Range("I61").Value = 8524.9375
$ws.Range("K61").Value = 8524.9375
$ws.Range("M61").Value = -8312.9375
$ws.Range("H63").Value = 4319.8
$ws.Range("J63").Value = 7549.5
$ws.Range("L63").Value = 7549.5
$ws.Range("N63").Value = -8921.5
$ws.Range("H66").Value = 4319.8
$ws.Range("J66").Value = 7549.5
$ws.Range("L66").Value = 37747.5
$ws.Range("N66").Value = -44611.5
$ws.Range("H80").Value = 39833.332
$ws.Range("J80").Value = 39750
$ws.Range("L80").Value = 39750
$ws.Range("N80").Value = -41746
$ws.Range("H83").Value = 39833.332
$ws.Range("J83").Value = 39750
$ws.Range("L83").Value = 119250
$ws.Range("N83").Value = -129234
$ws.Range("H88").Value = 2582.7144
$ws.Range("I88").Value = 1670.3334
$ws.Range("K88").Value = 1670.3334
$ws.Range("M88").Value = -1264.3334
$ws.Range("H91").Value = 2582.7144
$ws.Range("I91").Value = 1670.3334
$ws.Range("K91").Value = 1670.3334
$ws.Range("M91").Value = -266.3334
$ws.Range("H132").Value = 4557.1523
$ws.Range("I132").Value = 2962.4
$ws.Range("J132").Value = 6882.8335
$ws.Range("K132").Value = 8887.200000000001
$ws.Range("L132").Value = 20648.5005
$ws.Range("M132").Value = -6357.200000000001
$ws.Range("N132").Value = -25708.5005
$ws.Range("H136").Value = 7442.524
$ws.Range("I136").Value = 8524.9375
$ws.Range("K136").Value = 25574.8125
$ws.Range("M136").Value = -23024.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34832.668
$ws.Range("J35").Value = 34832.668
$ws.Range("L35").Value = 34832.668
$ws.Range("N35").Value = -35452.668
$ws.Range("H82").Value = 21108.7
$ws.Range("J82").Value = 46082.668
$ws.Range("L82").Value = 46082.668
$ws.Range("N82").Value = -46848.668
$ws.Range("H85").Value = 21108.7
$ws.Range("J85").Value = 46082.668
$ws.Range("L85").Value = 46082.668
$ws.Range("N85").Value = -48734.668
$ws.Range("H86").Value = 3622.818
$ws.Range("I86").Value = 3816.7778
$ws.Range("K86").Value = 3816.7778
$ws.Range("M86").Value = -2693.7778
$ws.Range("H89").Value = 3622.818
$ws.Range("I89").Value = 3816.7778
$ws.Range("K89").Value = 19083.889
$ws.Range("M89").Value = -13467.889
$ws.Range("H105").Value = 6260.607
$ws.Range("I105").Value = 3561.389
$ws.Range("K105").Value = 3561.389
$ws.Range("M105").Value = -1814.389
$ws.Range("H134").Value = 1309.3529
$ws.Range("I134").Value = 1190
$ws.Range("K134").Value = 3570
$ws.Range("M134").Value = -1035

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 626912.4
$ws.Range("I4").Value = 835250
$ws.Range("J4").Value = 1899.5
$ws.Range("K4").Value = 835250
$ws.Range("L4").Value = 1899.5
$ws.Range("M4").Value = -835138
$ws.Range("N4").Value = -2123.5
$ws.Range("H22").Value = 2406.3
$ws.Range("I22").Value = 817.7692
$ws.Range("J22").Value = 5356.4287
$ws.Range("K22").Value = 817.7692
$ws.Range("M22").Value = -467.7692
$ws.Range("N22").Value = -6056.4287
$ws.Range("H86").Value = 44213.54
$ws.Range("I86").Value = 121524
$ws.Range("K86").Value = 121524
$ws.Range("M86").Value = -120401
$ws.Range("H89").Value = 44213.54
$ws.Range("I89").Value = 121524
$ws.Range("K89").Value = 607620
$ws.Range("M89").Value = -602004
$ws.Range("H99").Value = 13782.947
$ws.Range("I99").Value = 25468.555
$ws.Range("J99").Value = 3265.9
$ws.Range("K99").Value = 25468.555
$ws.Range("L99").Value = 3265.9
$ws.Range("M99").Value = -23970.555
$ws.Range("N99").Value = -6261.9
$ws.Range("H107").Value = 814
$ws.Range("I107").Value = 720.05554
$ws.Range("J107").Value = 1055.5714
$ws.Range("K107").Value = 720.05554
$ws.Range("L107").Value = 1055.5714
$ws.Range("M107").Value = 1199.94446
$ws.Range("N107").Value = -4895.5714
$ws.Range("H126").Value = 13782.947
$ws.Range("I126").Value = 25468.555
$ws.Range("J126").Value = 3265.9
$ws.Range("K126").Value = 76405.66500000001
$ws.Range("L126").Value = 9797.700000000001
$ws.Range("M126").Value = -73935.66500000001
$ws.Range("N126").Value = -14737.7
$ws.Range("H132").Value = 4320.8433
$ws.Range("I132").Value = 5744.355
$ws.Range("J132").Value = 2114.4
$ws.Range("K132").Value = 17233.065
$ws.Range("L132").Value = 6343.200000000001
$ws.Range("M132").Value = -14703.065
$ws.Range("N132").Value = -11403.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4153934.5
$ws.Range("J4").Value = 3603.1428
$ws.Range("L4").Value = 10809.4284
$ws.Range("N4").Value = -11033.4284
$ws.Range("H38").Value = 486
$ws.Range("I38").Value = 545
$ws.Range("J38").Value = 250
$ws.Range("K38").Value = 1635
$ws.Range("L38").Value = 750
$ws.Range("M38").Value = -1288
$ws.Range("N38").Value = -1444
$ws.Range("H109").Value = 6423.625
$ws.Range("I109").Value = 2120.5
$ws.Range("K109").Value = 6361.5
$ws.Range("M109").Value = -5321.5
$ws.Range("H131").Value = 3325.4807
$ws.Range("I131").Value = 1655.4445
$ws.Range("J131").Value = 3675.0232
$ws.Range("K131").Value = 4966.333500000001
$ws.Range("L131").Value = 11025.0696
$ws.Range("M131").Value = 73.66649999999936
$ws.Range("N131").Value = -21105.0696

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("H132").Value = 9820.25
$ws.Range("I132").Value = 7473.65
$ws.Range("J132").Value = 15686.75
$ws.Range("K132").Value = 22420.95
$ws.Range("L132").Value = 47060.25
$ws.Range("M132").Value = -19890.95
$ws.Range("N132").Value = -52120.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3195.3713
$ws.Range("I46").Value = 1382.8334
$ws.Range("J46").Value = 5114.5293
$ws.Range("K46").Value = 1382.8334
$ws.Range("L46").Value = 5114.5293
$ws.Range("M46").Value = -1194.8334
$ws.Range("N46").Value = -5490.5293
$ws.Range("H55").Value = 1367.7354
$ws.Range("I55").Value = 1366.0834
$ws.Range("J55").Value = 1368.6364
$ws.Range("K55").Value = 1366.0834
$ws.Range("L55").Value = 1368.6364
$ws.Range("M55").Value = -1193.0834
$ws.Range("N55").Value = -1714.6364
$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H75").Value = 16407.334
$ws.Range("I75").Value = 16407.334
$ws.Range("K75").Value = 16407.334
$ws.Range("M75").Value = -15471.334
$ws.Range("H78").Value = 16407.334
$ws.Range("I78").Value = 16407.334
$ws.Range("K78").Value = 49222.00199999999
$ws.Range("M78").Value = -44542.00199999999
$ws.Range("H136").Value = 7357
$ws.Range("I136").Value = 5375
$ws.Range("K136").Value = 16125
$ws.Range("M136").Value = -13575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2175.8484
$ws.Range("I113").Value = 1715.8182
$ws.Range("J113").Value = 3095.9092
$ws.Range("K113").Value = 5147.4546
$ws.Range("L113").Value = 9287.7276
$ws.Range("M113").Value = -2977.4546
$ws.Range("N113").Value = -13627.7276
$ws.Range("H122").Value = 47159.16
$ws.Range("I122").Value = 1749.579
$ws.Range("J122").Value = 190956.17
$ws.Range("K122").Value = 5248.737
$ws.Range("L122").Value = 572868.51
$ws.Range("M122").Value = -2798.737
$ws.Range("N122").Value = -577768.51

